$wb = $excel.ActiveWorkbook

# --- Male_25m sheet: swap rows 7 and 8 (Albert Barrabino <-> Gabriel Rognes Steen) ---
$wsMale25 = $wb.Worksheets.Item("Male_25m")

$wsMale25.Cells.Item(7, 1).Value = "Gabriel Rognes Steen"
$wsMale25.Cells.Item(7, 2).Value = "1.00,69"
$wsMale25.Cells.Item(7, 4).Value = "28.09.2024"
$wsMale25.Cells.Item(7, 5).Value = "Bergen"

$wsMale25.Cells.Item(8, 1).Value = "Albert Barrabino"
$wsMale25.Cells.Item(8, 2).Value = "1.00,70"
$wsMale25.Cells.Item(8, 4).Value = "25.03.2012"
$wsMale25.Cells.Item(8, 5).Value = "Drammen"

# --- Female_25m sheet: update row 11 (Elise Lund -> Karoline Volden) ---
$wsFemale25 = $wb.Worksheets.Item("Female_25m")

$wsFemale25.Cells.Item(11, 1).Value = "Karoline Volden"
$wsFemale25.Cells.Item(11, 2).Value = "1.09,69"
$wsFemale25.Cells.Item(11, 4).Value = "30.09.2017"
$wsFemale25.Cells.Item(11, 5).Value = "Husebybadet"
